# Add a new row for LeetCode problem 2099:
# "Find Subsequence of Length K With the Largest Sum"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

$ws.Cells.Item($row, 1).Value = 2099
$ws.Cells.Item($row, 2).Value = "Find Subsequence of Length K With the Largest Sum"
$ws.Cells.Item($row, 3).Value = "#array #hash-table #sorting #heap "
$ws.Cells.Item($row, 4).Value = "easy"
$ws.Cells.Item($row, 5).Value = 1
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 13
$ws.Cells.Item($row, 8).Value = 45836
$ws.Cells.Item($row, 9).Value = 45836

# Copy styling from the row above so the new row matches the sheet formatting
$ws.Range("A19:I19").Copy()
$ws.Range("A20:I20").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A20").Select()
